# Add a new "Q8" column (J) to the sheet, matching the style used by the
# existing header cells (B1:I1), and populate the newly-available data
# points for rows 4 and 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell J1 = "Q8", using the same formatting as the other
# header cells (bold/centered/bordered style carried by e.g. I1).
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("J1").Value = "Q8"

# Row 4 gains values for columns G, H, I and the new column J.
$ws.Range("G4").Value = 0.2890697267702507
$ws.Range("H4").Value = -0.6507920071323952
$ws.Range("I4").Value = 0.4578003130087183
$ws.Range("J4").Value = -0.1119550751434417

# Row 8 gains values for columns G, H, I.
$ws.Range("G8").Value = 0.7010458975705092
$ws.Range("H8").Value = 0.6218889942996384
$ws.Range("I8").Value = 0.4230596606995932
